$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Quiz,Woksheet"
$ws.Range("C4").Value = "Quiz"
$ws.Range("C3").Value = "Worksheet"
$ws.Range("C1").Value = "Con"

$ws.Range("E4").Select()
